$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Corrected column A (Esfuerzo) values
$newValues = @(
    115.84686971533257,
    199.05705156897849,
    204.50283107021198,
    167.0955755316007,
    129.50778341584157,
    145.50147062310319,
    240.9018738756746,
    166.09720726453887,
    138.88165033783784,
    104.4924812030075,
    50.851287637727964,
    52.084008546600444,
    94.612820382424744,
    142.37455026777019,
    97.79465886454183,
    200.86981410069851
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}

# Column A loses its 2-decimal number format, matching column D's plain style
$ws.Range("A2:A17").NumberFormat = "General"

# Update selection to reflect the last-edited cell
$ws.Range("E15").Select()
